$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" note on Hoja1!A1 ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$nl = "`n"
$newText = "Conversión del día 💰" + $nl + `
"✅ Dólar paralelo: 68" + $nl + `
$nl + `
"Binance" + $nl + `
"✅ 1000 Bs = 5.95 = 24403.75 pesos" + $nl + `
"✅ 24403.75 pesos = 5.97 = 968.98 Bs" + $nl + `
$nl + `
"Promedio competencia" + $nl + `
"✅ Tasa pesos: 20" + $nl + `
"✅ Tasa Bs: 20" + $nl + `
"✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- Update the "tasas" sheet rate figures ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 168.091
$wsTasas.Range("O10").Value = 4102.05
$wsTasas.Range("O12").Value = 162.2
